$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column (C) for rows 2-9
# from serial date 45184 (2023-09-15) to 45185 (2023-09-16).
$ws.Range("C2:C9").Value = 45185
